$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new data row (row 8): a time value (12:00 noon => 0.5) in
#     column A and the re-used text "c" (same shared string as B4/B7) in
#     column B -------------------------------------------------------------
$ws.Range("A8").Value2 = 0.5
$ws.Range("B8").Value2 = "c"

# --- Re-extend the AutoFilter range to include the new row and restore the
#     filter criteria: the existing 0.046 / 0.516 plus the new 0.500
#     (formatted) value for the added row. This matches the fixed OOXML
#     import, which no longer double-converts the filter values (string ->
#     number -> string) and instead filters with the formatted cell value,
#     so "0.500" now shows up correctly instead of being dropped/garbled. --
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:B8").AutoFilter(1, @("0.046", "0.500", "0.516"), 7)

# Rows 2 and 6 remain filtered out/hidden, as before.
$ws.Rows(2).Hidden = $true
$ws.Rows(6).Hidden = $true

# Keep the defined "_FilterDatabase" name (driving the autofilter range) in
# sync with the grown range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Munka1!_FilterDatabase") {
        $n.RefersTo = "=Munka1!`$A`$1:`$B`$8"
    }
}

# --- Update the saved cursor/selection position (cosmetic, matches the
#     author's saved selection at the time of the edit) --------------------
[void]$ws.Range("C7").Select()
